$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.660.04'
$ws.Range("E2").Value = '  +4.60%  '
$ws.Range("D3").Value = '3.134.78'
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.41'
$ws.Range("E5").Value = '  +2.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '609.65'
$ws.Range("E6").Value = '  -0.62%  '
$ws.Range("E7").Value = '  +0.78%  '
$ws.Range("E8").Value = '  -1.72%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '3.134.74'
$ws.Range("E10").Value = '  +0.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.783'
$ws.Range("E11").Value = '  -6.65%  '
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("D13").Value = '97.325.01'
$ws.Range("E13").Value = '  +4.54%  '
$ws.Range("E14").Value = '  -2.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.43'
$ws.Range("E15").Value = '  +0.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '33.77'
$ws.Range("E16").Value = '  -3.93%  '
$ws.Range("D17").Value = '3.717.08'
$ws.Range("D18").Value = '3.132.25'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '518.09'
$ws.Range("E19").Value = '  +17.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.42'
$ws.Range("E20").Value = '  -9.95%  '
$ws.Range("E21").Value = '  -2.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.63'
$ws.Range("E22").Value = '  -6.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000191'
$ws.Range("E23").Value = '  -4.30%  '
$ws.Range("E24").Value = '  -4.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '88.24'
$ws.Range("E25").Value = '  +2.93%  '
$ws.Range("E26").Value = '  -4.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.52'
$ws.Range("E27").Value = '  -10.25%  '
$ws.Range("D28").Value = '3.299.04'
$ws.Range("E28").Value = '  +0.54%  '
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("E30").Value = '  -0.47%  '
$ws.Range("E31").Value = '  -3.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.122'
$ws.Range("E32").Value = '  -2.65%  '
$ws.Range("E33").Value = '  -0.76%  '
$ws.Range("E34").Value = '  -3.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '26.57'
$ws.Range("E35").Value = '  +2.51%  '
$ws.Range("E36").Value = '  -5.68%  '
$ws.Range("E37").Value = '  -9.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '24.33'
$ws.Range("E38").Value = '  +1.40%  '
$ws.Range("E39").Value = '  -1.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '468.59'
$ws.Range("E40").Value = '  -1.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.433'
$ws.Range("E41").Value = '  -2.43%  '
$ws.Range("E42").Value = '  -6.61%  '
$ws.Range("E43").Value = '  -10.36%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.09'
$ws.Range("E45").Value = '  -6.23%  '
$ws.Range("E46").Value = '  +1.92%  '
$ws.Range("E47").Value = '  -1.03%  '
$ws.Range("E48").Value = '  +1.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.45'
$ws.Range("E49").Value = '  +1.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '44.11'
$ws.Range("E50").Value = '  +0.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  +0.04%  '
